$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").Value = "Davide Rosà"
$ws.Range("B42").Value = "Elia Barozzi | I Magnifici"
$ws.Range("C42").Value = "Mattia Baldessarini | Shark Attack"
$ws.Range("D42").Value = "Michele Merighi | Clitoriders"
$ws.Range("E42").Value = "FEDERICO NICOLODI | U.S. Guarna"
$ws.Range("F42").Value = "Alberto Simoncelli | I Magnifici"
